$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.657985007885657
$ws.Range("C2").Value = 0.2409211227275136
$ws.Range("D2").Value = 0.4282344898476822
$ws.Range("E2").Value = 0.1010138446183664
$ws.Range("G2").Value = 0.002714825517674324
$ws.Range("I2").Value = 3.57611820949009
$ws.Range("J2").Value = 0.01813705359000117
$ws.Range("K2").Value = 2.592665901837677
$ws.Range("L2").Value = 0.7416679462038331
$ws.Range("N2").Value = 5.278814373002007

$ws.Range("B3").Value = 2.625694755306569
$ws.Range("C3").Value = 0.2331228688211695
$ws.Range("D3").Value = 0.4267626735653351
$ws.Range("E3").Value = 0.1011979672184449
$ws.Range("G3").Value = 0.00271998817399863
$ws.Range("I3").Value = 3.565383307327295
$ws.Range("J3").Value = 0.01801644599913654
$ws.Range("K3").Value = 2.553165932664172
$ws.Range("L3").Value = 0.7379445411949916
$ws.Range("N3").Value = 5.273354847626962

$ws.Range("B4").Value = 2.607306995407299
$ws.Range("C4").Value = 0.2284744758995032
$ws.Range("D4").Value = 0.4260403690709325
$ws.Range("E4").Value = 0.1013434587761548
$ws.Range("G4").Value = 0.002723325956178959
$ws.Range("I4").Value = 3.559834864318631
$ws.Range("J4").Value = 0.01794080490289041
$ws.Range("K4").Value = 2.530346814274907
$ws.Range("L4").Value = 0.7359832250159712
$ws.Range("N4").Value = 5.270854343902101

$ws.Range("B5").Value = 2.60017553168575
$ws.Range("C5").Value = 0.2266152121054859
$ws.Range("D5").Value = 0.4257916681294489
$ws.Range("E5").Value = 0.1014109197153186
$ws.Range("G5").Value = 0.002724728488098933
$ws.Range("I5").Value = 3.557835953975626
$ws.Range("J5").Value = 0.01790957762983236
$ws.Range("K5").Value = 2.521408162810786
$ws.Range("L5").Value = 0.735265677313393
$ws.Range("N5").Value = 5.270049282336402

$ws.Range("B6").Value = 2.599013202246027
$ws.Range("C6").Value = 0.2263085919509393
$ws.Range("D6").Value = 0.4257531292404622
$ws.Range("E6").Value = 0.101422615556368
$ws.Range("G6").Value = 0.002724963939817998
$ws.Range("I6").Value = 3.557519865333816
$ws.Range("J6").Value = 0.01790436792251526
$ws.Range("K6").Value = 2.51994566170535
$ws.Range("L6").Value = 0.7351514642869432
$ws.Range("N6").Value = 5.269928515311463

$ws.Range("B7").Value = 2.607209353554907
$ws.Range("C7").Value = 0.2284492597252665
$ws.Range("D7").Value = 0.4260368301509345
$ws.Range("E7").Value = 0.1013443354702286
$ws.Range("G7").Value = 0.002723344699633928
$ws.Range("I7").Value = 3.559806845134787
$ws.Range("J7").Value = 0.01794038539708698
$ws.Range("K7").Value = 2.530224805913065
$ws.Range("L7").Value = 0.7359732170778273
$ws.Range("N7").Value = 5.270842620790916

$ws.Range("B8").Value = 2.646552700645515
$ws.Range("C8").Value = 0.2382031816339918
$ws.Range("D8").Value = 0.4276893750788844
$ws.Range("E8").Value = 0.1010706059431339
$ws.Range("G8").Value = 0.00271657084106205
$ws.Range("I8").Value = 3.572200302803296
$ws.Range("J8").Value = 0.018095795328672
$ws.Range("K8").Value = 2.578748493531293
$ws.Range("L8").Value = 0.7403166938615584
$ws.Range("N8").Value = 5.276755010527538

$ws.Range("B9").Value = 2.735129890140342
$ws.Range("C9").Value = 0.2584477325623027
$ws.Range("D9").Value = 0.4323686804917344
$ws.Range("E9").Value = 0.1007906310708204
$ws.Range("G9").Value = 0.00270461301618922
$ws.Range("I9").Value = 3.604787163809334
$ws.Range("J9").Value = 0.01838813374931902
$ws.Range("K9").Value = 2.685303047832861
$ws.Range("L9").Value = 0.7514127628399336
$ws.Range("N9").Value = 5.29512073568344

$ws.Range("B10").Value = 2.807199174557184
$ws.Range("C10").Value = 0.2740156514618093
$ws.Range("D10").Value = 0.436683813306658
$ws.Range("E10").Value = 0.1007407976985224
$ws.Range("G10").Value = 0.002696626751648841
$ws.Range("I10").Value = 3.633797065650228
$ws.Range("J10").Value = 0.01859560069920452
$ws.Range("K10").Value = 2.770581183087927
$ws.Range("L10").Value = 0.7611402156928193
$ws.Range("N10").Value = 5.312765696774193

$ws.Range("B11").Value = 2.841510293941155
$ws.Range("C11").Value = 0.2812516498979107
$ws.Range("D11").Value = 0.4388375041981192
$ws.Range("E11").Value = 0.1007518349397749
$ws.Range("G11").Value = 0.002693165195407018
$ws.Range("I11").Value = 3.648099519443846
$ws.Range("J11").Value = 0.01868845095210681
$ws.Range("K11").Value = 2.810905126073862
$ws.Range("L11").Value = 0.7659084032043779
$ws.Range("N11").Value = 5.321699793326502

$ws.Range("B12").Value = 2.854722869436898
$ws.Range("C12").Value = 0.2840140982631283
$ws.Range("D12").Value = 0.4396804708453175
$ws.Range("E12").Value = 0.1007608491056953
$ws.Range("G12").Value = 0.002691878899800784
$ws.Range("I12").Value = 3.653674758594818
$ws.Range("J12").Value = 0.01872339553961844
$ws.Range("K12").Value = 2.826395450927521
$ws.Range("L12").Value = 0.7677633717486145
$ws.Range("N12").Value = 5.325213752936293

$ws.Range("B13").Value = 2.85186753035822
$ws.Range("C13").Value = 0.2834181594725464
$ws.Range("D13").Value = 0.4394977040050208
$ws.Range("E13").Value = 0.100758692918177
$ws.Range("G13").Value = 0.00269215483801441
$ws.Range("I13").Value = 3.652466946689245
$ws.Range("J13").Value = 0.01871587913297823
$ws.Range("K13").Value = 2.823049517288212
$ws.Range("L13").Value = 0.7673616754688055
$ws.Range("N13").Value = 5.324451135492666

$ws.Range("B14").Value = 2.842592896176825
$ws.Range("C14").Value = 0.2814784695804633
$ws.Range("D14").Value = 0.4389063062958343
$ws.Range("E14").Value = 0.100752479717924
$ws.Range("G14").Value = 0.002693058880619731
$ws.Range("I14").Value = 3.648555005667433
$ws.Range("J14").Value = 0.01869133016614999
$ws.Range("K14").Value = 2.812175102228309
$ws.Range("L14").Value = 0.7660600232249237
$ws.Range("N14").Value = 5.321986265275171

$ws.Range("B15").Value = 2.836940532472624
$ws.Range("C15").Value = 0.2802932672115048
$ws.Range("D15").Value = 0.4385476274133993
$ws.Range("E15").Value = 0.1007493031981461
$ws.Range("G15").Value = 0.002693615821878415
$ws.Range("I15").Value = 3.646179569159386
$ws.Range("J15").Value = 0.018676265228871
$ws.Range("K15").Value = 2.805542940976636
$ws.Range("L15").Value = 0.7652691523197745
$ws.Range("N15").Value = 5.320493506786789

$ws.Range("B16").Value = 2.804987594701288
$ws.Range("C16").Value = 0.27354587738688
$ws.Range("D16").Value = 0.4365469009873237
$ws.Range("E16").Value = 0.1007407535899176
$ws.Range("G16").Value = 0.002696856410855761
$ws.Range("I16").Value = 3.632884634093642
$ws.Range("J16").Value = 0.01858950230919465
$ws.Range("K16").Value = 2.767976753153619
$ws.Range("L16").Value = 0.7608355094420176
$ws.Range("N16").Value = 5.312200118977501

$ws.Range("B17").Value = 2.785776596314861
$ws.Range("C17").Value = 0.2694461594553275
$ws.Range("D17").Value = 0.4353683564119279
$ws.Range("E17").Value = 0.1007441328914158
$ws.Range("G17").Value = 0.002698888221972893
$ws.Range("I17").Value = 3.625011966622239
$ws.Range("J17").Value = 0.01853588758065694
$ws.Range("K17").Value = 2.745323416237227
$ws.Range("L17").Value = 0.7582035070370949
$ws.Range("N17").Value = 5.307345019172885

$ws.Range("B18").Value = 2.774870600999463
$ws.Range("C18").Value = 0.267102601545929
$ws.Range("D18").Value = 0.4347084405066113
$ws.Range("E18").Value = 0.1007492499354399
$ws.Range("G18").Value = 0.002700073010226634
$ws.Range("I18").Value = 3.620587877382107
$ws.Range("J18").Value = 0.01850490579843367
$ws.Range("K18").Value = 2.732437845648207
$ws.Range("L18").Value = 0.7567219436391639
$ws.Range("N18").Value = 5.304637874218315

$ws.Range("B19").Value = 2.77120268681972
$ws.Range("C19").Value = 0.266311596669766
$ws.Range("D19").Value = 0.4344880880061481
$ws.Range("E19").Value = 0.1007515278869704
$ws.Range("G19").Value = 0.002700476935702177
$ws.Range("I19").Value = 3.61910782303849
$ws.Range("J19").Value = 0.01849439104876272
$ws.Range("K19").Value = 2.728099736977015
$ws.Range("L19").Value = 0.7562258578587802
$ws.Range("N19").Value = 5.303735935284351

$ws.Range("B20").Value = 2.787806769975816
$ws.Range("C20").Value = 0.2698810803468348
$ws.Range("D20").Value = 0.4354919566617355
$ws.Range("E20").Value = 0.1007434448188196
$ws.Range("G20").Value = 0.002698670262400828
$ws.Range("I20").Value = 3.62583925415818
$ws.Range("J20").Value = 0.01854160983634046
$ws.Range("K20").Value = 2.747719993003955
$ws.Range("L20").Value = 0.7584803456752809
$ws.Range("N20").Value = 5.307853013952723

$ws.Range("B21").Value = 2.845311116807864
$ws.Range("C21").Value = 0.2820475959879332
$ws.Range("D21").Value = 0.4390792703248536
$ws.Range("E21").Value = 0.1007541735732804
$ws.Range("G21").Value = 0.002692792676922872
$ws.Range("I21").Value = 3.649699714914732
$ws.Range("J21").Value = 0.01869854661838222
$ws.Range("K21").Value = 2.815363193003691
$ws.Range("L21").Value = 0.766441010271663
$ws.Range("N21").Value = 5.322706704479884

$ws.Range("B22").Value = 2.884174211952541
$ws.Range("C22").Value = 0.2901293742334019
$ws.Range("D22").Value = 0.441583552148586
$ws.Range("E22").Value = 0.1007893596211566
$ws.Range("G22").Value = 0.002689094197824257
$ws.Range("I22").Value = 3.666222010101492
$ws.Range("J22").Value = 0.01879985791436667
$ws.Range("K22").Value = 2.860857528132897
$ws.Range("L22").Value = 0.7719314805511033
$ws.Range("N22").Value = 5.333177080541304

$ws.Range("B23").Value = 2.863314993963172
$ws.Range("C23").Value = 0.2858039992695751
$ws.Range("D23").Value = 0.4402323554751177
$ws.Range("E23").Value = 0.1007680062530749
$ws.Range("G23").Value = 0.002691055117107484
$ws.Range("I23").Value = 3.657318756082802
$ws.Range("J23").Value = 0.01874589980634589
$ws.Range("K23").Value = 2.836458562933899
$ws.Range("L23").Value = 0.768974780007639
$ws.Range("N23").Value = 5.327518948646002

$ws.Range("B24").Value = 2.78688849655731
$ws.Range("C24").Value = 0.2696844108108394
$ws.Range("D24").Value = 0.4354360220274174
$ws.Range("E24").Value = 0.1007437460083498
$ws.Range("G24").Value = 0.002698768750035944
$ws.Range("I24").Value = 3.62546491970916
$ws.Range("J24").Value = 0.01853902329671353
$ws.Range("K24").Value = 2.746636070409352
$ws.Range("L24").Value = 0.758355088469969
$ws.Range("N24").Value = 5.307623087498001

$ws.Range("B25").Value = 2.709941597636259
$ws.Range("C25").Value = 0.2528500923340005
$ws.Range("D25").Value = 0.4309487707764532
$ws.Range("E25").Value = 0.1008389511192433
$ws.Range("G25").Value = 0.002707706936239698
$ws.Range("I25").Value = 3.595082880107313
$ws.Range("J25").Value = 0.01831035157123928
$ws.Range("K25").Value = 2.655252428327174
$ws.Range("L25").Value = 0.7481346112598857
$ws.Range("N25").Value = 5.289424971490547
